# Reproduce the author's edit:
# - On the "sac à dos" sheet (Tableau3 table), remove the last data row
#   (the "Crédit" / "miam l'argent" row) from the table, which:
#     * clears the row's values (A7:G7)
#     * shrinks the table range from A1:G7 to A1:G6
#     * removes the now-unused shared strings ("Crédit", "miam l'argent")
#       and reindexes every other reference to shared strings automatically
# - Makes "sac à dos" the active sheet with F18 selected (instead of "equipement")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sac à dos")

# Remove the last row of data from the inventory table (Tableau3)
$ws.Range("A7:G7").ClearContents()

$lo = $ws.ListObjects.Item("Tableau3")
$lo.Resize($ws.Range("A1:G6"))

# The row that used to be the final row of the used range (row 18) no longer
# has any table-column cells backing it, so clear those cells entirely.
$ws.Range("A18:G18").Clear()

# Make the "sac à dos" sheet the active tab, with F18 selected
$ws.Activate()
$ws.Range("F18").Select()
